$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column: clone the header formatting from the neighboring
# "sum" header cell (G1) so H1 picks up the same bold/border/centered style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data row value for the new column (unstyled, like the other numeric cells)
$ws.Range("H2").Value = 0

$excel.CutCopyMode = 0
